$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: (New) -> (Nueva)
$ws.Range("D10").Value = "(Nueva)"

# Row 11: unchanged (D11 stays empty)

# Row 14: (New) -> (Nueva); E14 text changes
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("E14").Value = "A08:2021-Fallas en el Software y en la Integridad de los Datos"

# Row 15: E15 text changes (remove long suffix)
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"

# Row 16: (New) -> (Nueva)
$ws.Range("D16").Value = "(Nueva)"

# Row 17: E17 text changes
$ws.Range("E17").Value = "* A partir de la encuesta"

# Update the selection to match diff (E23 selected, though out of current data range)
$ws.Range("E23").Select()
